$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $range = $d.Content
    $found = $range.Find.Execute($find, $true, $true, $false, $false, $false, `
                                  $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Find failed for: $find"
    }
    $range.Text = $replace
}

# 1. Title change III -> II
Replace-Text "Laboratory of Materials Engineering III" "Laboratory of Materials Engineering II"

# 2. Ativacao date
Replace-Text "Ativação: 01/01/2020" "Ativação: 01/01/2025"

# 3. Objetivos paragraph
Replace-Text "Apresentar aos alunos, na prática, experimentos envolvendo tópicos de Pirometalurgia; Solidificação / Fundição e Eletrometalurgia." "Apresentar aos alunos, na prática, experimentos envolvendo tópicos de Solidificação/Fundição e Engenharia de Superfície."

# 4. Add new docente before Gilberto Carvalho Coelho (as its own run + <w:br/>,
#    matching the existing run-per-docente structure of the list).
#    Strategy: split the paragraph into three via InsertParagraphBefore (which
#    creates genuine new runs at the break points), insert the new text into
#    the middle paragraph, then re-join the paragraphs by deleting just the
#    paragraph-mark characters -- this preserves the run boundaries that
#    InsertParagraphBefore created, instead of Word merging the plain text
#    into the neighbouring run the way InsertBefore/Text= would.
$gilberto = $d.Content
$gilberto.Find.Execute("5009972 - Gilberto Carvalho Coelho") | Out-Null
$gilberto.Collapse(1)
$gilberto.Duplicate.InsertParagraphBefore() | Out-Null

$newRun = $d.Content
$newRun.Find.Execute("5009972 - Gilberto Carvalho Coelho") | Out-Null
$newRun.Collapse(1)
$newRun.InsertBefore("3586455 - Cassius Olivio Figueiredo Terra Ruchert`v")

$splitPoint = $d.Content
$splitPoint.Find.Execute("5009972 - Gilberto Carvalho Coelho") | Out-Null
$splitPoint.Collapse(1)
$splitPoint.Duplicate.InsertParagraphBefore() | Out-Null

# Re-find the paragraph that now holds only "Cassius..." and merge it with
# its neighbours by deleting the paragraph-mark characters at each boundary
# (highest offset first so earlier offsets stay valid).
$cassiusPara = $d.Content
$cassiusPara.Find.Execute("3586455 - Cassius Olivio Figueiredo Terra Ruchert") | Out-Null
$cassiusParaObj = $d.Paragraphs.Item(1)
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*Cassius*") {
        $cassiusParaObj = $d.Paragraphs.Item($i)
        break
    }
}
$afterMark = $d.Range($cassiusParaObj.Range.End - 1, $cassiusParaObj.Range.End)
$afterMark.Delete()
$carlosParaObj = $d.Paragraphs.Item($cassiusParaObj.Range.Start - 1)
# Recompute the paragraph that now ends right before "Cassius..." text
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*Carlos Angelo Nunes*") {
        $carlosParaObj = $d.Paragraphs.Item($i)
        break
    }
}
$beforeMark = $d.Range($carlosParaObj.Range.End - 1, $carlosParaObj.Range.End)
$beforeMark.Delete()

# 5. Aulas praticas paragraph
Replace-Text "Aulas práticas: pirometalurgia; eletrometalurgia e solidificação / fundição de metais e ligas." "Aulas práticas envolvendo experimentos de solidificação/fundição de metais e ligas e experimentos de modificação de superfície. Tipos de experimentos a serem realizados são exemplificados abaixo."

# 6. Big "Programa" paragraph replacement
$oldPrograma = "Eletrometalurgia:1) Eletrodeposição de cobre e/ou níquel;2) Anodização.PirometalurgiaOs experimentos poderão ser alterados a cada oferecimento da disciplina, mas envolverão normalmente os seguintes assuntos:1) Caracterização química e microestrutural de minérios; 2) Caracterização química e microestrutural de Ferro-Ligas;3) Secagem; calcinação de carbonatos e/ou hidróxidos; ustulação de sulfetos; aglomeração de minérios;4) Redução carbotérmica, por hidrogênio ou metalotérmica de óxidos; 5) Decomposição térmica sob vácuo de óxidos; 6) Refino de metais e/ou ligas.Solidificação / Fundição:Verificação experimental do efeito de algumas variáveis básicas de fundição tais como tipo de molde (cerâmicos, metálicos, ...), refrigeração ou pré-aquecimento do molde, temperatura de vazamento, adição de agentes inoculantes, agitação mecânica, dentre outras, nas micro e macroestruturas de metais e ligas metálicas. A turma será dividida em grupos sendo que cada um dos grupos deverá investigar com mais detalhe uma das variáveis (ou combinação de variáveis) do processamento por fundição."
$newPrograma = "1.Solidificação / Fundição:Solidificação dos metais. Processos de montagem, tecnologia de fundição. Segregação e defeitos em peças fundidas. Qualidade de peças fundidas. Verificação experimental do efeito de algumas variáveis básicas de fundição tais como tipo de molde (cerâmicos, metálicos, ...), refrigeração ou pré-aquecimento do molde, temperatura de vazamento, adição de agentes inoculantes, agitação mecânica, dentre outras, nas micro e macroestruturas de metais e ligas metálicas. 2.Experimentos de modificação superficial tais como: eletrodeposição, nitretação, carbonitretação, cementação, deposição por soldagem, aspersão Térmica, PVD, CVD, pack cementation, têmpera superficial, ensaios tribológicos (abrasão a dois ou mais corpos e roda de borracha), caracterizações microestruturais.A turma será dividida em grupos, cada um recebendo orientação para realização de tarefas específicas. Dependendo do tipo de experimento, os resultados de todos os grupos serão compartilhados e cada grupo deverá desenvolver sua interpretação sobre o conjunto daqueles resultados.Cada grupo deverá também desenvolver um trabalho (não experimental), em tema livre, procurando integrar temáticas de fundição/solidificação e modificação superficial."
Replace-Text $oldPrograma $newPrograma

# 7. Metodo value
Replace-Text "Experimentos desenvolvidos em laboratórios, realização de relatórios para cada experimento e apresentação dos resultados obtidos." "Cada grupo de alunos receberá uma nota baseada em: 1) documento apresentando os experimentos desenvolvidos, resultados obtidos e conclusões; 2) apresentação oral dos experimentos desenvolvidos, resultados obtidos e conclusões; 3) documento relativo ao trabalho de integração; 4) apresentação oral do trabalho de integração."

# 8. Criterio value
Replace-Text "A avaliação será feita através das notas dos relatórios escritos e das apresentações orais dos trabalhos desenvolvidos em cada módulo da disciplina." "A avaliação será feita através das notas atribuídas aos documentos produzidos e às apresentações orais realizadas"

# 9. Norma de recuperacao - add comma
Replace-Text "Devido às características da disciplina não será oferecida recuperação" "Devido às características da disciplina, não será oferecida recuperação"

# 10. Bibliografia paragraph replacement
$oldBib = 'DENARO, A.R. Fundamentos de Eletroquímica. São Paulo: Editora Edgard Blucher, 1974. PLETCHER, D.; WALSH, F.C. Industrial Electrochemistry, 2ª Ed. Springer, 1990.KUHN, A.T. Industrial Electrochemical Processes. Elsevier Pub., 1971.ROSENQVIST, T. Principles of Extrative Metallurgy, McGraw-Hill, 1983.HABASHI, F. Extractive Metallurgy, Gordon and Breach Science Publishers, 1986.GARCIA, A. Solidificação: Fundamentos e Aplicações, Editora da Unicamp, 2001. FLEMINGS, M.C. Solidification Processing, McGraw-Hill, 1974.'
$newBib = 'GARCIA, A. Solidificação: Fundamentos e Aplicações, Editora da Unicamp, 2001.FLEMINGS, M.C. Solidification Processing, McGraw-Hill, 1974.CAMPOS FILHO, M.P. de, DAVIES, G.J. - Solidificação e fundição de metais e suas ligas. S.Paulo, EDUSP, 1978.SIEGEL, M. - Fundição. S.Paulo, ABM, 1979.Krauss, G. "Steels: Heat Treating and Processing Principles" - ASM International, 1990.ASM Handbook, vol. 5 – Surface Engineering. ISBN: 978-0-87170-384-2PLETCHER, D.; WALSH, F.C. Industrial Electrochemistry, 2ª Ed. Springer, 1990.'
Replace-Text $oldBib $newBib
